# Applies the "Changes of 29th march 2022" edit to the FedExShipments test
# fixture:
#  - refresh the ShipmentTracking (P), ActualRate (Q) and Result (R) columns
#    for rows 2-25 with the new batch of tracking numbers / PASS-FAIL outcomes
#  - move the saved view (scroll position + selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be stored as text (Excel's normal type-inference would
    # otherwise turn an all-digit string into a number, or a "$"-prefixed
    # string into a currency value), then drop back to the default "Normal"
    # style so no stray number-format style sticks to the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("P2") '320018165080'
$ws.Range("R2").Value = "PASS"

# Row 3
Set-TextValue $ws.Range("P3") '320018164874'
$ws.Range("R3").Value = "PASS"

# Row 4
Set-TextValue $ws.Range("P4") '320018164922'
$ws.Range("R4").Value = "PASS"

# Row 5
Set-TextValue $ws.Range("P5") '320018164944'
$ws.Range("R5").Value = "PASS"

# Row 6
Set-TextValue $ws.Range("P6") '320018164988'
$ws.Range("R6").Value = "PASS"

# Row 7
Set-TextValue $ws.Range("P7") '320018165002'
$ws.Range("R7").Value = "FAIL"

# Row 8
Set-TextValue $ws.Range("P8") '320018164360'
$ws.Range("R8").Value = "PASS"

# Row 9
Set-TextValue $ws.Range("P9") '320018164381'
$ws.Range("R9").Value = "PASS"

# Row 10
Set-TextValue $ws.Range("P10") '320018164418'
$ws.Range("R10").Value = "PASS"

# Row 11
Set-TextValue $ws.Range("P11") '320018164451'
$ws.Range("R11").Value = "PASS"

# Row 12
Set-TextValue $ws.Range("P12") '320018164495'
$ws.Range("R12").Value = "PASS"

# Row 13
Set-TextValue $ws.Range("P13") '320018164510'
$ws.Range("R13").Value = "PASS"

# Row 14
Set-TextValue $ws.Range("P14") '320018164543'
$ws.Range("R14").Value = "PASS"

# Row 15
Set-TextValue $ws.Range("P15") '320018164565'
$ws.Range("R15").Value = "PASS"

# Row 16
Set-TextValue $ws.Range("P16") '320018164598'
$ws.Range("R16").Value = "PASS"

# Row 17
Set-TextValue $ws.Range("P17") '320018164613'
$ws.Range("R17").Value = "PASS"

# Row 18
Set-TextValue $ws.Range("P18") '320018164657'
$ws.Range("R18").Value = "FAIL"

# Row 19
Set-TextValue $ws.Range("P19") '320018164679'
$ws.Range("R19").Value = "PASS"

# Row 20
Set-TextValue $ws.Range("P20") '320018164705'
$ws.Range("R20").Value = "FAIL"

# Row 21
Set-TextValue $ws.Range("P21") '320018164727'
$ws.Range("R21").Value = "PASS"

# Row 22 (also updates ActualRate)
Set-TextValue $ws.Range("P22") '320018164750'
Set-TextValue $ws.Range("Q22") '$195.48'
$ws.Range("R22").Value = "FAIL"

# Row 23
Set-TextValue $ws.Range("P23") '320018164760'
$ws.Range("R23").Value = "PASS"

# Row 24
Set-TextValue $ws.Range("P24") '320018164771'
$ws.Range("R24").Value = "FAIL"

# Row 25
Set-TextValue $ws.Range("P25") '320018164782'
$ws.Range("R25").Value = "PASS"

# Update the sheet view: column J becomes the left-most visible column, and
# the live selection moves to a single cell, S8.
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("S8").Select()
